$d = $word.ActiveDocument

function Find-ParagraphByText($doc, [string]$text) {
    $target = $null
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -eq ($text + [char]13)) {
            $target = $p
            break
        }
    }
    return $target
}

# =========================================================================
# Insertion point 1: after the "Phasor analysis" paragraph, add two new
# paragraphs:
#   - "IQ DC transceiver"                                    (plain)
#   - "IQ mismatch, EVM, image, calibration"                 (bulleted list)
# =========================================================================

$anchor1 = Find-ParagraphByText $d "Phasor analysis"
$idx1 = $anchor1.Index

$r1 = $anchor1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$r1.InsertParagraphAfter()

$para1 = $d.Paragraphs($idx1 + 1)
$para1.Range.Text = "IQ DC transceiver"

$para2 = $d.Paragraphs($idx1 + 2)
$para2.Range.Text = "IQ mismatch, EVM, image, calibration"

# Apply a real bullet-list template (this also creates word/numbering.xml
# with a full hybridMultilevel bullet definition).
$gallery = $word.ListGalleries.Item(1)
$template = $gallery.ListTemplates.Item(1)
$para2.Range.ListFormat.ApplyListTemplateWithLevel($template, $false, 2, $false, $false)

# Attach the "ListParagraph" paragraph style cleanly (without disturbing the
# numPr that was just set) by re-inserting the paragraph's OOXML directly.
$listXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>IQ mismatch, EVM, image, calibration</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$d.Paragraphs($idx1 + 2).Range.Text = ""
$d.Paragraphs($idx1 + 2).Range.InsertXML($listXml) | Out-Null

# =========================================================================
# Insertion point 2: after the "Orthogonality of tones" paragraph, add
# three new plain paragraphs.
# =========================================================================

$anchor2 = Find-ParagraphByText $d "Orthogonality of tones"
$idx2 = $anchor2.Index

$r2 = $anchor2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r2.InsertParagraphAfter()
$r2.InsertParagraphAfter()

$d.Paragraphs($idx2 + 1).Range.Text = "Geometric series"
$d.Paragraphs($idx2 + 2).Range.Text = "Trigonometric identities"
$d.Paragraphs($idx2 + 3).Range.Text = "Stochastic signal processing " + [char]0x2013 + " does noise power drop as you average more samples? Noise variance should drop (CLT)."

Write-Host "done"
